# Insert a new data row at row 62 (shifts existing rows 62-178 down to 63-179)
# and populate it with the new record (same values as the old row 62 except
# for the date in column D, which becomes 2021-09-30 / serial 44469).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(62).Insert()

$ws.Range("A62").Value = 7
$ws.Range("B62").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C62").Value = 'Ñuble'
$ws.Range("D62").Value = 44469
$ws.Range("E62").Value = 16
$ws.Range("F62").Value = 100114013
$ws.Range("G62").Value = 'Zanahoria'
$ws.Range("H62").Value = 'Sin especificar'
$ws.Range("I62").Value = 'Primera'
$ws.Range("J62").Value = 120
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 6500
$ws.Range("M62").Value = 6250
$ws.Range("N62").Value = '$/saco 20 kilos'
$ws.Range("O62").Value = 'Provincia de Diguillín'
$ws.Range("P62").Value = 312
$ws.Range("Q62").Value = 20
$ws.Range("R62").Value = 'Hortaliza'
